# Append a new data row (ID_Boiler = 2, type "gases") to the
# OperationScenario_Component_Boiler input table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "gases"
$ws.Range("C3").Value = 15000
$ws.Range("D3").Value = "W"
$ws.Range("E3").Value = 0.35
$ws.Range("F3").Value = 35
$ws.Range("G3").Value = 55

# Leave the selection where it ended up after the edit.
$ws.Range("E7").Select()
